# Auto-generated update of leve profit calculations (Excalibur server leve data)
# Updates derived from refreshed Universalis market price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1776.4036
$ws.Range("J17").Value = 1886.5625
$ws.Range("L17").Value = 5659.6875
$ws.Range("N17").Value = -5995.6875

$ws.Range("H28").Value = 534.0625
$ws.Range("I28").Value = 436.33334
$ws.Range("K28").Value = 436.33334
$ws.Range("M28").Value = 48.66665999999998

$ws.Range("H76").Value = 7093.75
$ws.Range("I76").Value = 4500
$ws.Range("K76").Value = 4500
$ws.Range("M76").Value = -4185

$ws.Range("H79").Value = 7093.75
$ws.Range("I79").Value = 4500
$ws.Range("K79").Value = 4500
$ws.Range("M79").Value = -3408

$ws.Range("H80").Value = 5209174
$ws.Range("J80").Value = 1182.4
$ws.Range("L80").Value = 3547.2
$ws.Range("N80").Value = -5543.200000000001

$ws.Range("H83").Value = 5209174
$ws.Range("J83").Value = 1182.4
$ws.Range("L83").Value = 10641.6
$ws.Range("N83").Value = -20625.6

$ws.Range("H112").Value = 1762.8
$ws.Range("J112").Value = 1839
$ws.Range("L112").Value = 5517
$ws.Range("N112").Value = -7733

$ws.Range("H132").Value = 1513.1091
$ws.Range("I132").Value = 1263.8654
$ws.Range("K132").Value = 3791.5962
$ws.Range("M132").Value = -1261.5962

$ws.Range("H137").Value = 1573100
$ws.Range("I137").Value = 1723.25
$ws.Range("K137").Value = 5169.75
$ws.Range("M137").Value = -2619.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2863.48
$ws.Range("I32").Value = 1489.5172
$ws.Range("J32").Value = 12058.462
$ws.Range("K32").Value = 1489.5172
$ws.Range("L32").Value = 12058.462
$ws.Range("M32").Value = -1202.5172
$ws.Range("N32").Value = -12632.462

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 598
$ws.Range("I94").Value = 607.8
$ws.Range("J94").Value = 500
$ws.Range("K94").Value = 607.8
$ws.Range("L94").Value = 500
$ws.Range("M94").Value = -156.8
$ws.Range("N94").Value = -1402

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2762.3333
$ws.Range("I16").Value = 2715
$ws.Range("J16").Value = 2999
$ws.Range("K16").Value = 2715
$ws.Range("L16").Value = 2999
$ws.Range("M16").Value = -2428
$ws.Range("N16").Value = -3573

$ws.Range("H113").Value = 2762.3333
$ws.Range("I113").Value = 2715
$ws.Range("J113").Value = 2999
$ws.Range("K113").Value = 2715
$ws.Range("L113").Value = 2999
$ws.Range("M113").Value = -545
$ws.Range("N113").Value = -7339

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 237.25
$ws.Range("I2").Value = 131.6
$ws.Range("J2").Value = 272.46667
$ws.Range("K2").Value = 131.6
$ws.Range("L2").Value = 272.46667
$ws.Range("M2").Value = -18.59999999999999
$ws.Range("N2").Value = -498.46667

$ws.Range("H11").Value = 344368.88
$ws.Range("I11").Value = 12348.692
$ws.Range("J11").Value = 2502500
$ws.Range("K11").Value = 12348.692
$ws.Range("L11").Value = 2502500
$ws.Range("M11").Value = -12209.692
$ws.Range("N11").Value = -2502778

$ws.Range("H25").Value = 15000
$ws.Range("J25").Value = 15000
$ws.Range("L25").Value = 15000
$ws.Range("N25").Value = -16058

$ws.Range("H70").Value = 5401.5
$ws.Range("I70").Value = 5351.8887
$ws.Range("J70").Value = 5550.3335
$ws.Range("K70").Value = 5351.8887
$ws.Range("L70").Value = 5550.3335
$ws.Range("M70").Value = -5081.8887
$ws.Range("N70").Value = -6090.3335

$ws.Range("H73").Value = 5401.5
$ws.Range("I73").Value = 5351.8887
$ws.Range("J73").Value = 5550.3335
$ws.Range("K73").Value = 5351.8887
$ws.Range("L73").Value = 5550.3335
$ws.Range("M73").Value = -4415.8887
$ws.Range("N73").Value = -7422.3335

$ws.Range("H113").Value = 3289.9333
$ws.Range("I113").Value = 1557.4
$ws.Range("J113").Value = 4156.2
$ws.Range("K113").Value = 1557.4
$ws.Range("L113").Value = 4156.2
$ws.Range("M113").Value = 612.5999999999999
$ws.Range("N113").Value = -8496.200000000001

$ws.Range("H122").Value = 4792.636
$ws.Range("I122").Value = 5220
$ws.Range("J122").Value = 2086
$ws.Range("K122").Value = 15660
$ws.Range("L122").Value = 6258
$ws.Range("M122").Value = -13210
$ws.Range("N122").Value = -11158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 89999.75
$ws.Range("I4").Value = 67499.5
$ws.Range("J4").Value = 112500
$ws.Range("K4").Value = 67499.5
$ws.Range("L4").Value = 112500
$ws.Range("M4").Value = -67386.5
$ws.Range("N4").Value = -112726

$ws.Range("H11").Value = 5153.6665
$ws.Range("I11").Value = 5006
$ws.Range("J11").Value = 5227.5
$ws.Range("K11").Value = 5006
$ws.Range("L11").Value = 5227.5
$ws.Range("M11").Value = -4866
$ws.Range("N11").Value = -5507.5

$ws.Range("H25").Value = 52500
$ws.Range("I25").Value = 55000
$ws.Range("J25").Value = 50000
$ws.Range("K25").Value = 55000
$ws.Range("L25").Value = 50000
$ws.Range("M25").Value = -54770
$ws.Range("N25").Value = -50460

$ws.Range("H28").Value = 89999.75
$ws.Range("I28").Value = 67499.5
$ws.Range("J28").Value = 112500
$ws.Range("K28").Value = 67499.5
$ws.Range("L28").Value = 112500
$ws.Range("M28").Value = -67267.5
$ws.Range("N28").Value = -112964

$ws.Range("H37").Value = 89999.75
$ws.Range("I37").Value = 67499.5
$ws.Range("J37").Value = 112500
$ws.Range("K37").Value = 67499.5
$ws.Range("L37").Value = 112500
$ws.Range("M37").Value = -67392.5
$ws.Range("N37").Value = -112714

$ws.Range("H40").Value = 4411.5713
$ws.Range("J40").Value = 4295.6665
$ws.Range("L40").Value = 4295.6665
$ws.Range("N40").Value = -4567.6665

$ws.Range("H82").Value = 1978.15
$ws.Range("I82").Value = 1332.3
$ws.Range("K82").Value = 1332.3
$ws.Range("M82").Value = -971.3

$ws.Range("H85").Value = 1978.15
$ws.Range("I85").Value = 1332.3
$ws.Range("K85").Value = 1332.3
$ws.Range("M85").Value = -84.29999999999995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7569.75
$ws.Range("I62").Value = 6777.5
$ws.Range("K62").Value = 6777.5
$ws.Range("M62").Value = -6153.5

$ws.Range("H65").Value = 7569.75
$ws.Range("I65").Value = 6777.5
$ws.Range("K65").Value = 33887.5
$ws.Range("M65").Value = -30767.5

$ws.Range("H81").Value = 740.8570999999999
$ws.Range("I81").Value = 669.2
$ws.Range("K81").Value = 1338.4
$ws.Range("M81").Value = -277.4000000000001

$ws.Range("H84").Value = 740.8570999999999
$ws.Range("I84").Value = 669.2
$ws.Range("K84").Value = 6692
$ws.Range("M84").Value = -1388
